# "Se procesan de nuevo los datos con las nuevas dimensiones curadas"
#
# Column F ("provincia") was previously treated as an sdmx dimension
# (sdmx-dimension:refArea / dim / URI-Provincia). It is now re-processed
# as a curated iaest measure, matching the pattern used by the other
# measure columns (iaest-measure:<name> / medida / xsd:int).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "iaest-measure:provincia"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"
